# adding lab8 guide and IPAM update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate the new IPAM rows (38-44) in the exact cell-entry sequence the
# author used, so the shared-string table is rebuilt with the same
# first-seen ordering as the authored workbook.
$ws.Cells.Item(38,1).Value = "H5"
$ws.Cells.Item(38,3).Value = "2.2.2.2/30"
$ws.Cells.Item(38,4).Value = "2222::2/64"

$ws.Cells.Item(39,1).Value = "R6"

$ws.Cells.Item(40,2).Value = "eth4"

$ws.Cells.Item(42,1).Value = "R7"

$ws.Cells.Item(43,3).Value = "2.2.2.1/30"
$ws.Cells.Item(43,4).Value = "2222::1/64"

$ws.Cells.Item(39,4).Value = "2666::1/64"
$ws.Cells.Item(40,4).Value = "2666::2/64"
$ws.Cells.Item(41,4).Value = "2777::1/64"
$ws.Cells.Item(42,4).Value = "2777::2/64"

$ws.Cells.Item(40,3).Value = "10.40.100.1/24"
$ws.Cells.Item(42,3).Value = "10.40.101.2/24"
$ws.Cells.Item(39,3).Value = "10.40.100.2/24"
$ws.Cells.Item(41,3).Value = "10.40.101.1/24"

$ws.Cells.Item(44,1).Value = "R6 "
$ws.Cells.Item(44,3).Value = "10.40.6.1/32"

# Remaining cells reuse already-known strings (Interface names / hostnames).
$ws.Cells.Item(38,2).Value = "eth1"
$ws.Cells.Item(39,2).Value = "eth1"
$ws.Cells.Item(40,1).Value = "R4"
$ws.Cells.Item(41,1).Value = "R6"
$ws.Cells.Item(41,2).Value = "eth2"
$ws.Cells.Item(42,2).Value = "eth1"
$ws.Cells.Item(43,1).Value = "R7"
$ws.Cells.Item(43,2).Value = "eth2"
$ws.Cells.Item(44,2).Value = "Loopback"

# Update the sheet view to reflect where the author ended up scrolled/selected.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 19
$ws.Range("E22").Select()
